# deep sea double count fix
# Updates landings-area status table values after correcting a double
# count of deep sea catch. Also fixes the number format on row 20
# (Deep Sea) columns C:F to match the 3-decimal style used elsewhere
# (style index 7 / "#,##0.000") instead of the previous 2-decimal style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("C4").Value = 0.9859867053471901
$ws.Range("D4").Value = 0.2254840312804062
$ws.Range("E4").Value = 1.201838525915919
$ws.Range("F4").Value = 0.2254840312804062
$ws.Range("G4").Value = 15.12284798418123
$ws.Range("H4").Value = 69.0794593258551
$ws.Range("I4").Value = 15.79769268996366
$ws.Range("J4").Value = 84.20230731003633
$ws.Range("K4").Value = 15.79769268996366

# --- Row 5 ---
$ws.Range("D5").Value = 0.9839436444806937
$ws.Range("F5").Value = 0.9839436444806937
$ws.Range("G5").Value = 4.053027175415057
$ws.Range("H5").Value = 82.61771444705293
$ws.Range("I5").Value = 13.32925837753201
$ws.Range("J5").Value = 86.67074162246799
$ws.Range("K5").Value = 13.32925837753201

# --- Row 20 (Deep Sea) ---
# C20:F20 switch number format from 2-decimal to 3-decimal style.
$ws.Range("C20:F20").NumberFormat = "#,##0.000"
$ws.Range("C20").Value = 0.04298276364640884
$ws.Range("D20").Value = 0.03504234798678382
$ws.Range("E20").Value = 0.04298276364640884
$ws.Range("F20").Value = 0.03504234798678382
$ws.Range("H20").Value = 55.08837186735091
$ws.Range("I20").Value = 44.91162813264909
$ws.Range("J20").Value = 55.08837186735091
$ws.Range("K20").Value = 44.91162813264909

# --- Row 22 (Sharks) ---
$ws.Range("B22").Value = 0.02678414
$ws.Range("C22").Value = 0.02084913
$ws.Range("D22").Value = 0.00601671
$ws.Range("E22").Value = 0.04763327000000001
$ws.Range("F22").Value = 0.00601671
$ws.Range("G22").Value = 49.92385831271513
$ws.Range("H22").Value = 38.86139379735091
$ws.Range("I22").Value = 11.21474788993398
$ws.Range("J22").Value = 88.78525211006604
$ws.Range("K22").Value = 11.21474788993398

# --- Row 23 (Global) ---
$ws.Range("B23").Value = 18.50865755204853
$ws.Range("C23").Value = 33.60113412567404
$ws.Range("D23").Value = 17.34828565888902
$ws.Range("E23").Value = 52.10979167772258
$ws.Range("F23").Value = 17.34828565888902
$ws.Range("G23").Value = 26.6472356589297
$ws.Range("H23").Value = 48.37613624522654
$ws.Range("I23").Value = 24.97662809584375
$ws.Range("J23").Value = 75.02337190415624
$ws.Range("K23").Value = 24.97662809584375
